# LogBook.xlsx update — Joao/Duarte logbook entry for day 31/03 (44286)
$wb = $excel.ActiveWorkbook

# --- TimeTable_Team1: fill in the new work-session row ---
$wsTeam1 = $wb.Worksheets.Item("TimeTable_Team1")

$wsTeam1.Range("A2").Value = 44286
$wsTeam1.Range("B2").Value = 0.58333333333333337
$wsTeam1.Range("C2").Value = "Estudo do seguidor de linha. Uso do DMA para ler 2 inputs do sensor"
$wsTeam1.Range("D2").Value = 0.80555555555555547

# Make TimeTable_Team1 the active/selected sheet with B6 selected
# (this also clears the tabSelected flag previously on TimeTable_Team3)
$wsTeam1.Select()
$wsTeam1.Range("B6").Select()
